$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (e.g. A1) to the new headers
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the record values for every data row (rows 2 through 53)
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 111   # AD
    $ws.Cells.Item($r, 31).Value = 51    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
